$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows 21-27 down to 22-28
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly data point
$ws.Cells.Item(21, 1).Value = 10
$ws.Cells.Item(21, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(21, 3).Value = "La Araucanía"
$ws.Cells.Item(21, 4).Value = 44755
$ws.Cells.Item(21, 4).NumberFormat = $ws.Cells.Item(22, 4).NumberFormat
$ws.Cells.Item(21, 5).Value = 9
$ws.Cells.Item(21, 6).Value = 100112042
$ws.Cells.Item(21, 7).Value = "Locoto"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 90
$ws.Cells.Item(21, 11).Value = 3300
$ws.Cells.Item(21, 12).Value = 3300
$ws.Cells.Item(21, 13).Value = 3300
$ws.Cells.Item(21, 14).Value = "$/kilo"
$ws.Cells.Item(21, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 16).Value = 3300
$ws.Cells.Item(21, 17).Value = 1
$ws.Cells.Item(21, 18).Value = "Hortaliza"
